$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.200.09'
$ws.Range("E2").Value = '  +8.07%  '
$ws.Range("D3").Value = '1.590.98'
$ws.Range("E3").Value = '  +8.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9906'
$ws.Range("E5").Value = '  +3.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '298.28'
$ws.Range("E6").Value = '  +7.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3627'
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3347'
$ws.Range("E8").Value = '  +8.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.27'
$ws.Range("E9").Value = '  +4.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.118'
$ws.Range("E10").Value = '  +4.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06948'
$ws.Range("E11").Value = '  +4.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.42'
$ws.Range("E13").Value = '  +6.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.822'
$ws.Range("E14").Value = '  +5.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.536'
$ws.Range("E15").Value = '  +5.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9911'
$ws.Range("E16").Value = '  +3.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001062'
$ws.Range("E17").Value = '  +3.45%  '
$ws.Range("D18").Value = '1.590.85'
$ws.Range("E18").Value = '  +7.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06560'
$ws.Range("E19").Value = '  +10.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '76.03'
$ws.Range("E20").Value = '  +10.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.83'
$ws.Range("E21").Value = '  +8.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.918'
$ws.Range("E22").Value = '  +7.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.62'
$ws.Range("E23").Value = '  +4.14%  '
$ws.Range("D24").Value = '22.186.59'
$ws.Range("E24").Value = '  +7.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.379'
$ws.Range("E25").Value = '  +5.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.500'
$ws.Range("E26").Value = '  +17.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.02'
$ws.Range("E27").Value = '  +3.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.15'
$ws.Range("E28").Value = '  +11.74%  '
$ws.Range("D29").Value = '1.754.75'
$ws.Range("E29").Value = '  +7.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.26'
$ws.Range("E30").Value = '  +6.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.981'
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.891'
$ws.Range("E32").Value = '  +18.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9201'
$ws.Range("E33").Value = '  +14.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08135'
$ws.Range("E34").Value = '  +1.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.615'
$ws.Range("E35").Value = '  +6.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.68'
$ws.Range("E36").Value = '  +12.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.112'
$ws.Range("E37").Value = '  +8.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.238'
$ws.Range("E38").Value = '  +1.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.349'
$ws.Range("E39").Value = '  +12.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06002'
$ws.Range("E40").Value = '  +4.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02177'
$ws.Range("E41").Value = '  +5.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1983'
$ws.Range("E42").Value = '  +6.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9921'
$ws.Range("E43").Value = '  +3.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5782'
$ws.Range("E44").Value = '  +9.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.766'
$ws.Range("E45").Value = '  +7.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.84'
$ws.Range("E46").Value = '  +5.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.49'
$ws.Range("E47").Value = '  +5.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5556'
$ws.Range("E48").Value = '  +6.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.939'
$ws.Range("E49").Value = '  +6.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06713'
$ws.Range("E50").Value = '  +3.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.46'
$ws.Range("E51").Value = '  +8.21%  '
